$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 8.937933333333332
$ws.Range("H2").Value = 26.8138
$ws.Range("I2").Value = 0.2302024600837126
$ws.Range("J2").Value = 0.2302024600837126
$ws.Range("M2").Value = 10.67287833333334
$ws.Range("N2").Value = 32.018635
$ws.Range("O2").Value = 0.2083050184412124
$ws.Range("P2").Value = 0.2083050184412124
$ws.Range("Q2").Value = 95.3934750181111
$ws.Range("R2").Value = 858.541275163
$ws.Range("S2").Value = 0.04795232769295021
$ws.Range("T2").Value = 0.0479523276929502
$ws.Range("G3").Value = 8.937933333333332
$ws.Range("H3").Value = 26.8138
$ws.Range("I3").Value = 0.2302024600837126
$ws.Range("J3").Value = 0.2302024600837126
$ws.Range("O3").Value = 0.5074067008595954
$ws.Range("P3").Value = 0.5074067008595954
$ws.Range("Q3").Value = 232.3673659169777
$ws.Range("R3").Value = 2091.3062932528
$ws.Range("S3").Value = 0.1168062708008393
$ws.Range("T3").Value = 0.1168062708008393
$ws.Range("G4").Value = 8.937933333333332
$ws.Range("H4").Value = 26.8138
$ws.Range("I4").Value = 0.2302024600837126
$ws.Range("J4").Value = 0.2302024600837126
$ws.Range("M4").Value = 14.566016
$ws.Range("N4").Value = 43.698048
$ws.Range("O4").Value = 0.2842882806991923
$ws.Range("P4").Value = 0.2842882806991923
$ws.Range("Q4").Value = 130.1900799402666
$ws.Range("R4").Value = 1171.7107194624
$ws.Range("S4").Value = 0.06544386158992309
$ws.Range("T4").Value = 0.06544386158992309
$ws.Range("I5").Value = 0.5278886986241245
$ws.Range("J5").Value = 0.5278886986241244
$ws.Range("M5").Value = 10.67287833333334
$ws.Range("N5").Value = 32.018635
$ws.Range("O5").Value = 0.2083050184412124
$ws.Range("P5").Value = 0.2083050184412124
$ws.Range("Q5").Value = 218.7515171046884
$ws.Range("R5").Value = 1968.763653942195
$ws.Range("S5").Value = 0.1099618651018058
$ws.Range("T5").Value = 0.1099618651018058
$ws.Range("I6").Value = 0.5278886986241245
$ws.Range("J6").Value = 0.5278886986241244
$ws.Range("O6").Value = 0.5074067008595954
$ws.Range("P6").Value = 0.5074067008595954
$ws.Range("S6").Value = 0.2678542629899322
$ws.Range("T6").Value = 0.2678542629899321
$ws.Range("I7").Value = 0.5278886986241245
$ws.Range("J7").Value = 0.5278886986241244
$ws.Range("M7").Value = 14.566016
$ws.Range("N7").Value = 43.698048
$ws.Range("O7").Value = 0.2842882806991923
$ws.Range("P7").Value = 0.2842882806991923
$ws.Range("Q7").Value = 298.545340690304
$ws.Range("R7").Value = 2686.908066212736
$ws.Range("S7").Value = 0.1500725705323864
$ws.Range("T7").Value = 0.1500725705323864
$ws.Range("G8").Value = 9.392449999999998
$ws.Range("H8").Value = 28.17735
$ws.Range("I8").Value = 0.241908841292163
$ws.Range("J8").Value = 0.2419088412921629
$ws.Range("M8").Value = 10.67287833333334
$ws.Range("N8").Value = 32.018635
$ws.Range("O8").Value = 0.2083050184412124
$ws.Range("P8").Value = 0.2083050184412124
$ws.Range("Q8").Value = 100.2444761019167
$ws.Range("R8").Value = 902.20028491725
$ws.Range("S8").Value = 0.05039082564645633
$ws.Range("T8").Value = 0.05039082564645631
$ws.Range("G9").Value = 9.392449999999998
$ws.Range("H9").Value = 28.17735
$ws.Range("I9").Value = 0.241908841292163
$ws.Range("J9").Value = 0.2419088412921629
$ws.Range("O9").Value = 0.5074067008595954
$ws.Range("P9").Value = 0.5074067008595954
$ws.Range("Q9").Value = 244.1838380990666
$ws.Range("R9").Value = 2197.6545428916
$ws.Range("S9").Value = 0.1227461670688239
$ws.Range("T9").Value = 0.1227461670688239
$ws.Range("G10").Value = 9.392449999999998
$ws.Range("H10").Value = 28.17735
$ws.Range("I10").Value = 0.241908841292163
$ws.Range("J10").Value = 0.2419088412921629
$ws.Range("M10").Value = 14.566016
$ws.Range("N10").Value = 43.698048
$ws.Range("O10").Value = 0.2842882806991923
$ws.Range("P10").Value = 0.2842882806991923
$ws.Range("Q10").Value = 136.8105769792
$ws.Range("R10").Value = 1231.2951928128
$ws.Range("S10").Value = 0.06877184857688277
$ws.Range("T10").Value = 0.06877184857688277
